$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("downlink")
$ws.Columns("E:E").Insert()
$ws.Rows("8:8").Insert()
$ws.Range("A7:L7").Copy()
$ws.Range("A8:L8").PasteSpecial(-4122)
$ws.Rows("26:26").Insert()
Write-Host ("C25: " + $ws.Range("C25").Value2)
Write-Host ("C26: " + $ws.Range("C26").Value2)
Write-Host ("C27: " + $ws.Range("C27").Value2)
